# Updated TPM-derived values for the Fgf2-Sdc1 ligand-receptor pair sheet.
# Each row corresponds to a Sending-cluster/Target-cluster combination;
# ligand (G-J), receptor (M-P) and edge-weight (Q-T) statistics are refreshed
# with the new TPM-based figures, leaving all other cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.7321483333333333
$ws.Range("H2").Value = 2.196445
$ws.Range("I2").Value = 0.05113520435363902
$ws.Range("J2").Value = 0.05113520435363902
$ws.Range("M2").Value = 1.815761
$ws.Range("N2").Value = 5.447283000000001
$ws.Range("O2").Value = 0.07007596730428067
$ws.Range("P2").Value = 0.07007596730428067
$ws.Range("Q2").Value = 1.329406389881667
$ws.Range("R2").Value = 11.964657508935
$ws.Range("S2").Value = 0.003583348908383318
$ws.Range("T2").Value = 0.003583348908383318

# Row 3
$ws.Range("G3").Value = 0.7321483333333333
$ws.Range("H3").Value = 2.196445
$ws.Range("I3").Value = 0.05113520435363902
$ws.Range("J3").Value = 0.05113520435363902
$ws.Range("O3").Value = 0.5079540516959071
$ws.Range("P3").Value = 0.5079540516959072
$ws.Range("Q3").Value = 9.63636162393111
$ws.Range("R3").Value = 86.72725461537999
$ws.Range("S3").Value = 0.02597433423572913
$ws.Range("T3").Value = 0.02597433423572914

# Row 4
$ws.Range("G4").Value = 0.7321483333333333
$ws.Range("H4").Value = 2.196445
$ws.Range("I4").Value = 0.05113520435363902
$ws.Range("J4").Value = 0.05113520435363902
$ws.Range("M4").Value = 9.711409333333334
$ws.Range("N4").Value = 29.134228
$ws.Range("O4").Value = 0.3747940411327002
$ws.Range("P4").Value = 0.3747940411327002
$ws.Range("Q4").Value = 7.110192157717778
$ws.Range("R4").Value = 63.99172941945999
$ws.Range("S4").Value = 0.01916516988384681
$ws.Range("T4").Value = 0.01916516988384681

# Row 5
$ws.Range("G5").Value = 0.7321483333333333
$ws.Range("H5").Value = 2.196445
$ws.Range("I5").Value = 0.05113520435363902
$ws.Range("J5").Value = 0.05113520435363902
$ws.Range("M5").Value = 1.222391
$ws.Range("N5").Value = 3.667173
$ws.Range("O5").Value = 0.04717593986711188
$ws.Range("P5").Value = 0.04717593986711189
$ws.Range("Q5").Value = 0.8949715333316666
$ws.Range("R5").Value = 8.054743799984999
$ws.Range("S5").Value = 0.002412351325679752
$ws.Range("T5").Value = 0.002412351325679752

# Row 6
$ws.Range("I6").Value = 0.7165747117895102
$ws.Range("J6").Value = 0.7165747117895102
$ws.Range("M6").Value = 1.815761
$ws.Range("N6").Value = 5.447283000000001
$ws.Range("O6").Value = 0.07007596730428067
$ws.Range("P6").Value = 0.07007596730428067
$ws.Range("Q6").Value = 18.62941612773267
$ws.Range("R6").Value = 167.664745149594
$ws.Range("S6").Value = 0.05021466607443606
$ws.Range("T6").Value = 0.05021466607443606

# Row 7
$ws.Range("I7").Value = 0.7165747117895102
$ws.Range("J7").Value = 0.7165747117895102
$ws.Range("O7").Value = 0.5079540516959071
$ws.Range("P7").Value = 0.5079540516959072
$ws.Range("S7").Value = 0.3639870281963086
$ws.Range("T7").Value = 0.3639870281963087

# Row 8
$ws.Range("I8").Value = 0.7165747117895102
$ws.Range("J8").Value = 0.7165747117895102
$ws.Range("M8").Value = 9.711409333333334
$ws.Range("N8").Value = 29.134228
$ws.Range("O8").Value = 0.3747940411327002
$ws.Range("P8").Value = 0.3747940411327002
$ws.Range("Q8").Value = 99.63749946023378
$ws.Range("R8").Value = 896.737495142104
$ws.Range("S8").Value = 0.2685679320050904
$ws.Range("T8").Value = 0.2685679320050904

# Row 9
$ws.Range("I9").Value = 0.7165747117895102
$ws.Range("J9").Value = 0.7165747117895102
$ws.Range("M9").Value = 1.222391
$ws.Range("N9").Value = 3.667173
$ws.Range("O9").Value = 0.04717593986711188
$ws.Range("P9").Value = 0.04717593986711189
$ws.Range("Q9").Value = 12.54153526251267
$ws.Range("R9").Value = 112.873817362614
$ws.Range("S9").Value = 0.03380508551367496
$ws.Range("T9").Value = 0.03380508551367496

# Row 10
$ws.Range("G10").Value = 2.568000333333333
$ws.Range("H10").Value = 7.704001
$ws.Range("I10").Value = 0.1793560346266988
$ws.Range("J10").Value = 0.1793560346266988
$ws.Range("M10").Value = 1.815761
$ws.Range("N10").Value = 5.447283000000001
$ws.Range("O10").Value = 0.07007596730428067
$ws.Range("P10").Value = 0.07007596730428067
$ws.Range("Q10").Value = 4.662874853253667
$ws.Range("R10").Value = 41.96587367928301
$ws.Range("S10").Value = 0.01256854761832597
$ws.Range("T10").Value = 0.01256854761832597

# Row 11
$ws.Range("G11").Value = 2.568000333333333
$ws.Range("H11").Value = 7.704001
$ws.Range("I11").Value = 0.1793560346266988
$ws.Range("J11").Value = 0.1793560346266988
$ws.Range("O11").Value = 0.5079540516959071
$ws.Range("P11").Value = 0.5079540516959072
$ws.Range("Q11").Value = 33.79940749125378
$ws.Range("R11").Value = 304.194667421284
$ws.Range("S11").Value = 0.09110462448474306
$ws.Range("T11").Value = 0.09110462448474307

# Row 12
$ws.Range("G12").Value = 2.568000333333333
$ws.Range("H12").Value = 7.704001
$ws.Range("I12").Value = 0.1793560346266988
$ws.Range("J12").Value = 0.1793560346266988
$ws.Range("M12").Value = 9.711409333333334
$ws.Range("N12").Value = 29.134228
$ws.Range("O12").Value = 0.3747940411327002
$ws.Range("P12").Value = 0.3747940411327002
$ws.Range("Q12").Value = 24.93890240513645
$ws.Range("R12").Value = 224.450121646228
$ws.Range("S12").Value = 0.06722157301927693
$ws.Range("T12").Value = 0.06722157301927693

# Row 13
$ws.Range("G13").Value = 2.568000333333333
$ws.Range("H13").Value = 7.704001
$ws.Range("I13").Value = 0.1793560346266988
$ws.Range("J13").Value = 0.1793560346266988
$ws.Range("M13").Value = 1.222391
$ws.Range("N13").Value = 3.667173
$ws.Range("O13").Value = 0.04717593986711188
$ws.Range("P13").Value = 0.04717593986711189
$ws.Range("Q13").Value = 3.139100495463667
$ws.Range("R13").Value = 28.251904459173
$ws.Range("S13").Value = 0.008461289504352779
$ws.Range("T13").Value = 0.008461289504352779

# Row 14
$ws.Range("G14").Value = 0.7579039999999999
$ws.Range("H14").Value = 2.273712
$ws.Range("I14").Value = 0.05293404923015203
$ws.Range("J14").Value = 0.05293404923015203
$ws.Range("M14").Value = 1.815761
$ws.Range("N14").Value = 5.447283000000001
$ws.Range("O14").Value = 0.07007596730428067
$ws.Range("P14").Value = 0.07007596730428067
$ws.Range("Q14").Value = 1.376172524944
$ws.Range("R14").Value = 12.385552724496
$ws.Range("S14").Value = 0.003709404703135317
$ws.Range("T14").Value = 0.003709404703135317

# Row 15
$ws.Range("G15").Value = 0.7579039999999999
$ws.Range("H15").Value = 2.273712
$ws.Range("I15").Value = 0.05293404923015203
$ws.Range("J15").Value = 0.05293404923015203
$ws.Range("O15").Value = 0.5079540516959071
$ws.Range("P15").Value = 0.5079540516959072
$ws.Range("Q15").Value = 9.975351561578664
$ws.Range("R15").Value = 89.77816405420799
$ws.Range("S15").Value = 0.02688806477912634
$ws.Range("T15").Value = 0.02688806477912634

# Row 16
$ws.Range("G16").Value = 0.7579039999999999
$ws.Range("H16").Value = 2.273712
$ws.Range("I16").Value = 0.05293404923015203
$ws.Range("J16").Value = 0.05293404923015203
$ws.Range("M16").Value = 9.711409333333334
$ws.Range("N16").Value = 29.134228
$ws.Range("O16").Value = 0.3747940411327002
$ws.Range("P16").Value = 0.3747940411327002
$ws.Range("Q16").Value = 7.360315979370666
$ws.Range("R16").Value = 66.242843814336
$ws.Range("S16").Value = 0.01983936622448597
$ws.Range("T16").Value = 0.01983936622448597

# Row 17
$ws.Range("G17").Value = 0.7579039999999999
$ws.Range("H17").Value = 2.273712
$ws.Range("I17").Value = 0.05293404923015203
$ws.Range("J17").Value = 0.05293404923015203
$ws.Range("M17").Value = 1.222391
$ws.Range("N17").Value = 3.667173
$ws.Range("O17").Value = 0.04717593986711188
$ws.Range("P17").Value = 0.04717593986711189
$ws.Range("Q17").Value = 0.9264550284639999
$ws.Range("R17").Value = 8.338095256175999
$ws.Range("S17").Value = 0.002497213523404392
$ws.Range("T17").Value = 0.002497213523404392
